# Deploy the implementation guide: refresh generated metadata sheet.
$wb = $excel.ActiveWorkbook

# 1. Rename the second sheet ("Include from Ferlab.bio CodeS" -> "Include #0")
$wsInclude = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$wsInclude.Name = "Include #0"

# 2. Update the Metadata sheet values
$ws = $wb.Worksheets.Item("Metadata")

# Regenerated "Date" value (row 8)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Regenerated "Contact" value (row 10) - now resolves to a display string
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# A "Jurisdiction" property is now emitted right after "Contact", pushing
# Description/Purpose/Copyright/Immutable down by one row each.
$ws.Rows.Item(11).Insert()

# Carry the standard body-row formatting onto the newly inserted row.
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
